# Quizvragen aangepast via Admin
# Adds a new "image_url" column (L) to every sheet of the quiz workbook,
# fills in an image URL for the first question on the "Wiskunde 3" sheet,
# and updates that question's text/choices.

$wb = $excel.ActiveWorkbook

function Add-ImageUrlColumn {
    param(
        $ws,
        [int]$lastRow
    )

    # New header cell L1: copy the formatting of K1 (bold header style) then
    # set its text so the style (s="1") is preserved exactly.
    $ws.Cells.Item(1, 11).Copy($ws.Cells.Item(1, 12))
    $ws.Cells.Item(1, 12).Value = "image_url"

    # Data rows: create an (empty) cell in column L for every existing row
    # so the column is fully populated, matching the other columns.
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 12).NumberFormat = "@"
        $ws.Cells.Item($r, 12).Value = ""
        $ws.Cells.Item($r, 12).Style = "Normal"
    }
}

# --- Sheet "DC" (5 data rows incl. header) ---
$wsDC = $wb.Worksheets.Item("DC")
Add-ImageUrlColumn $wsDC 5

# --- Sheet "Wiskunde 3" (21 data rows incl. header) ---
$wsWiskunde = $wb.Worksheets.Item("Wiskunde 3")
Add-ImageUrlColumn $wsWiskunde 21

# Question q1 on "Wiskunde 3" was reworked from a goniometry question into
# an Ohm's law question, and now references an image.
$wsWiskunde.Cells.Item(2, 4).Value = "Wat is de juiste formule voor de stroom I?"
$wsWiskunde.Cells.Item(2, 5).Value = "['I = U/R', ' U = I*R', ' R = U/I']"
$wsWiskunde.Cells.Item(2, 12).Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/Wiskunde_3_q0_1763116823.jpg"

# --- Sheet "AC" (2 data rows incl. header) ---
$wsAC = $wb.Worksheets.Item("AC")
Add-ImageUrlColumn $wsAC 2
